$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that should now show a checkmark ("ü" in Wingdings renders as a tick)
# in column C: C14, C20, C30.
foreach ($addr in @("C14", "C20", "C30")) {
    $cell = $ws.Range($addr)
    $cell.Value = "ü"
    $cell.Font.Name = "Wingdings"
    $cell.Font.Size = 12
}

# C28 previously had the checkmark; it should now be blank again, using the
# same plain font as the other unchecked cells in column C.
$cell = $ws.Range("C28")
$cell.Value = ""
$cell.Font.Name = "Calibri"
$cell.Font.Size = 12
